$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2xxx905"
$ws.Range("B2").Value = "bil******hotmail.com"
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 5.5

$ws.Range("A3").Value = "2xxx105"
$ws.Range("B3").Value = "877******qq.com"
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 6.5

$ws.Range("A4").Value = "2xxx858"
$ws.Range("B4").Value = "she**************163.com"
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 7.5

$ws.Range("A5").Value = "2xxx917"
$ws.Range("B5").Value = "378********qq.com"
$ws.Range("C5").Value = 24.56
$ws.Range("D5").Value = 8

$ws.Range("A6").Value = "2xxx556"
$ws.Range("B6").Value = "Zfa*****************gmail.com"
$ws.Range("C6").Value = 23.1
$ws.Range("D6").Value = 2.5

$ws.Range("A7").Value = "2xxx485"
$ws.Range("B7").Value = "290*******qq.com"
$ws.Range("C7").Value = 21.3
$ws.Range("D7").Value = 4.6

$ws.Range("A8").Value = "2xxx647"
$ws.Range("B8").Value = "skt******outlook.com"
$ws.Range("C8").Value = 19.6
$ws.Range("D8").Value = 3.6

$ws.Range("B9").Value = "784********qq.com"
$ws.Range("A9").Value = "2xxx957"
$ws.Range("C9").Value = 18.9
$ws.Range("D9").Value = 4.2

$ws.Range("A10").Value = "2xxx441"
$ws.Range("B10").Value = "233*******qq.com"
$ws.Range("C10").Value = 15.86
$ws.Range("D10").Value = 6.6

$ws.Range("A11").Value = "2xxx845"
$ws.Range("B11").Value = "fei*****163.com"
$ws.Range("C11").Value = 14.79
$ws.Range("D11").Value = 2

$ws.Range("E2").Formula = "=C2*0.9+D2*0.01"
$ws.Range("E3").Formula = "=C3*0.9+D3*0.01"
$ws.Range("E4:E11").Formula = "=C4*0.9+D4*0.01"

$ws.Range("B11").ClearFormats()
